$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook-level metadata ---

# --- Clear stray duplicate formatting on existing rows 82-101 (col A/B) ---
$ws.Range("A82:A101").ClearFormats()
$ws.Range("B84:B101").ClearFormats()
# Re-apply Times New Roman to col D so the engine folds it onto the surviving style slot
$ws.Range("D82:D101").Font.Name = "Times New Roman"

# --- Append new GamePlay requirement rows (102-116) ---
$ws.Range("A102").Value = 99
$ws.Range("B102").Value = 36
$ws.Range("C102").Value = 'Game Play shall determine whether a move is legal or not and decides whether a piece has been captured.'
$ws.Range("D102").Value = 'GamePlay.GamePlay.returnBoard'
$ws.Rows.Item(102).RowHeight = 28.5

$ws.Range("A103").Value = 100
$ws.Range("B103").Value = 36
$ws.Range("C103").Value = 'A checker''s piece shall not move to a spot taken by another checker''s piece.  If the user attempts to move in a spot that is taken send an error message to Game Manager'
$ws.Range("D103").Value = 'GamePlay.rules.isLegal'
$ws.Rows.Item(103).RowHeight = 42.75

$ws.Range("A104").Value = 101
$ws.Range("B104").Value = 36
$ws.Range("C104").Value = 'A checker''s piece shall only move diagonally up the board unless it is marked by a crown, allowing the piece to move diagonally up and down.'
$ws.Range("D104").Value = 'GamePlay.rules.isLegal'
$ws.Rows.Item(104).RowHeight = 42.75

$ws.Range("A105").Value = 102
$ws.Range("B105").Value = 36
$ws.Range("C105").Value = 'A checker''s piece shall not move beyond the border of the playing field.'
$ws.Range("D105").Value = 'GamePlay.rules.inBounds'
$ws.Rows.Item(105).RowHeight = 28.5

$ws.Range("A106").Value = 103
$ws.Range("B106").Value = 36
$ws.Range("C106").Value = 'Game Play shall communicate with Game Manager, notifying whether the requested move is legal and whether a piece has been captured.'
$ws.Range("D106").Value = 'GamePlay.GamePlay.returnBoard'
$ws.Rows.Item(106).RowHeight = 42.75

$ws.Range("A107").Value = 104
$ws.Range("B107").Value = 36
$ws.Range("C107").Value = 'A player shall not be allowed to move checker pieces of their oponent'
$ws.Range("D107").Value = 'GamePlay.rules.canMovePiece'
$ws.Rows.Item(107).RowHeight = 28.5

$ws.Range("A108").Value = 105
$ws.Range("B108").Value = 36
$ws.Range("C108").Value = 'The user is able to make multiple successive jumps during one turn, if each jump captures an opponents piece and follows the rules associated with the user''s cheker piece (i.e a king could do this in either a backwards or forwards direction, and a regular piece is only able to do this in a forward direction)'
$ws.Range("D108").Value = 'GamePlay.rules.isLegal'
$ws.Rows.Item(108).RowHeight = 85.5

$ws.Range("A109").Value = 106
$ws.Range("B109").Value = 36
$ws.Range("C109").Value = 'The user must capture an opponent''s piece if a legal move is presented/available'
$ws.Range("D109").Value = 'GamePlay.rules.isLegal / GamePlay.GamePlay.returnBoard'
$ws.Rows.Item(109).RowHeight = 28.5

$ws.Range("A110").Value = 107
$ws.Range("B110").Value = 36
$ws.Range("C110").Value = 'The bot should be subject to the same rules as the user'
$ws.Range("D110").Value = 'GamePlay.GamePlay.returnBoard'

$ws.Range("A111").Value = 108
$ws.Range("B111").Value = 36
$ws.Range("C111").Value = 'A king checker piece is not able to move both up and down during one turn. It must choose to either move up or down'
$ws.Range("D111").Value = 'GamePlay.rules.isLegal'
$ws.Rows.Item(111).RowHeight = 28.5

$ws.Range("A112").Value = 109
$ws.Range("B112").Value = 36
$ws.Range("C112").Value = 'If the user is presented two different opportunities to capture an opponents checker piece during their turn, they are allowed to chose which checker''s piece they capture'
$ws.Range("D112").Value = 'GamePlay.rules.isLegal'
$ws.Rows.Item(112).RowHeight = 57

$ws.Range("A113").Value = 110
$ws.Range("B113").Value = 36
$ws.Range("C113").Value = 'A king checker piece shall be held to the same capture rules as a regular checker piece. '
$ws.Range("D113").Value = 'GamePlay.rules.isLegal'
$ws.Rows.Item(113).RowHeight = 28.5

$ws.Range("A114").Value = 111
$ws.Range("B114").Value = 36
$ws.Range("C114").Value = 'A captured piece shall be removed from the game'
$ws.Range("D114").Value = 'GamePlay.rules.isLegal'

$ws.Range("A115").Value = 112
$ws.Range("B115").Value = 36
$ws.Range("C115").Value = 'A pieces move length shall be proportional to the amount of the opponent''s pieces the moving piece jumps over.'
$ws.Range("D115").Value = 'GamePlay.rules.returnBoard'
$ws.Rows.Item(115).RowHeight = 28.5

$ws.Range("A116").Value = 113
$ws.Range("B116").Value = 36
$ws.Range("C116").Value = 'GamePlay shall return a list of possible moves a player can take during their turn'
$ws.Range("D116").Value = 'GamePlay.GamePlay.returnMoves'
$ws.Rows.Item(116).RowHeight = 28.5

# --- Update selection / view to reflect the new bottom of the sheet ---
$ws.Range("D117").Select()
